$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.352.54'
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.896.55'
$ws.Range("E3").Value = '  +2.28%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.87'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.05'
$ws.Range("E6").Value = '  +1.93%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.896.60'
$ws.Range("E7").Value = '  +2.40%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -1.44%  '

$ws.Range("E10").Value = '  -1.96%  '

$ws.Range("E11").Value = '  +1.83%  '

$ws.Range("E12").Value = '  -0.31%  '

$ws.Range("E13").Value = '  +3.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.32'
$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.546.79'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.899.25'
$ws.Range("E16").Value = '  +1.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.469.88'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.47'
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.32'
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.09'
$ws.Range("E21").Value = '  -2.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.49'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.726'
$ws.Range("E23").Value = '  +0.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000166'
$ws.Range("E24").Value = '  +3.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.67'
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.01'
$ws.Range("E27").Value = '  -2.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("E28").Value = '  +1.26%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  -1.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.046.05'
$ws.Range("E31").Value = '  +2.19%  '

$ws.Range("E32").Value = '  -1.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.73'
$ws.Range("E33").Value = '  -3.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.77'
$ws.Range("E34").Value = '  -0.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.847.14'
$ws.Range("E35").Value = '  +2.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.107'
$ws.Range("E36").Value = '  -0.53%  '

$ws.Range("E37").Value = '  +1.11%  '

$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("E40").Value = '  +3.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("E42").Value = '  -1.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '431.68'
$ws.Range("E43").Value = '  +0.46%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.98'
$ws.Range("E44").Value = '  -0.57%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.13'
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.53'
$ws.Range("E46").Value = '  +1.75%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.83'
$ws.Range("E48").Value = '  +1.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000270'
$ws.Range("E49").Value = '  +17.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.800.46'
$ws.Range("E50").Value = '  -1.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.43'
$ws.Range("E51").Value = '  -0.31%  '
